$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.375.94'
$ws.Range("E2").Value = '  +0.22%  '
$ws.Range("D3").Value = '1.836.93'
$ws.Range("E3").Value = '  -0.47%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9996'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -1.47%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '243.37'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.40%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6244'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.81%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.001'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -1.02%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07397'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.66%  '
$ws.Range("E9").Value = '  -0.80%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '23.28'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.11%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07646'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.06%  '
$ws.Range("D12").Value = '1.834.70'
$ws.Range("E12").Value = '  -0.11%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.013'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.22%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6756'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.33%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '83.05'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.24%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.000009355'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +3.29%  '
$ws.Range("E17").Value = '  -0.45%  '
$ws.Range("D18").Value = '29.356.04'
$ws.Range("E18").Value = '  +0.35%  '
$ws.Range("D19").Value = '2.093.49'
$ws.Range("E19").Value = '  +0.97%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '238.17'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.40%  '
$ws.Range("E21").Value = '  -1.46%  '
$ws.Range("E22").Value = '  -1.10%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.377'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.60%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.002'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.84%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '158.62'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.98%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1412'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.82%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.463'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.85%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '17.74'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.14%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.06077'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +8.47%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.494'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.92%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.251'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.69%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.114'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.26%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.093'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.53%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.858'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.14%  '
$ws.Range("E35").Value = '  -0.29%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7251'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -3.86%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.613'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.02%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.894'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.08%  '
$ws.Range("D39").Value = '1.219.55'
$ws.Range("E39").Value = '  +0.46%  '
$ws.Range("E40").Value = '  -1.40%  '
$ws.Range("E41").Value = '  -2.21%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.9113'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.32%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.002'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.95%  '
$ws.Range("D44").Value = '1.999.41'
$ws.Range("E44").Value = '  +1.34%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '101.73'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.62%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '65.49'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.20%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.00000000120'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -3.04%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.5069'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.37%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.254'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.60%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.4053'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.34%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.1147'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +3.18%  '
